$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 (shifts existing rows 11-22 down to 12-23)
# and populate it with the new weekly record.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44483
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112013
$ws.Range("G11").Value = "Alcachofa"
$ws.Range("H11").Value = "Madrigal"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = "$/caja 40 unidades"
$ws.Range("O11").Value = "Región de Coquimbo"
$ws.Range("P11").Value = 362
$ws.Range("Q11").Value = 40
$ws.Range("R11").Value = "Hortaliza"
